$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.741.50"
$ws.Range("E2").Value = "  +0.23%  "
$ws.Range("D3").Value = "2.563.00"
$ws.Range("E3").Value = "  +1.22%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.81"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.60%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "101.45"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.55%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.570"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.31%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.530"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.85%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.25"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +2.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0805"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.40"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -0.68%  "
$ws.Range("E13").Value = "  +0.72%  "
$ws.Range("D14").Value = "2.954.68"
$ws.Range("E14").Value = "  +1.30%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "16.33"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +9.20%  "
$ws.Range("D16").Value = "2.597.81"
$ws.Range("E16").Value = "  +3.87%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.844"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.19%  "
$ws.Range("D18").Value = "42.766.77"
$ws.Range("E18").Value = "  +0.30%  "
$ws.Range("E19").Value = "  -0.13%  "
$ws.Range("B20").Value = "InternetComputer(DFINITY)"
$ws.Range("C20").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.37"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.97%  "
$ws.Range("B21").Value = "ShibaInu"
$ws.Range("C21").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D21").Value = "0.0₃0956"
$ws.Range("E21").Value = "  +0.14%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "243.25"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -3.58%  "
$ws.Range("E24").Value = "  -0.80%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.07"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +0.95%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "26.48"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.19%  "
$ws.Range("E27").Value = "  +0.04%  "
$ws.Range("B28").Value = "Toncoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.36"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.60%  "
$ws.Range("B29").Value = "InjectiveProtocol"
$ws.Range("C29").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "40.13"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.20%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "10.14"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.52%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "159.27"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +2.00%  "
$ws.Range("E32").Value = "  -1.88%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "2.78"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +14.57%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0804"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +1.19%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.07"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.66%  "
$ws.Range("E36").Value = "  -2.99%  "
$ws.Range("E37").Value = "  -2.10%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "18.21"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.57%  "
$ws.Range("E39").Value = "  +0.19%  "
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "4.23"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +11.61%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.86"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.26%  "
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.34"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +3.80%  "
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +0.15%  "
$ws.Range("E45").Value = "  -0.57%  "
$ws.Range("D46").Value = "1.964.74"
$ws.Range("E46").Value = "  -0.68%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.00"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +0.53%  "
$ws.Range("D48").Value = "2.805.91"
$ws.Range("E48").Value = "  +1.20%  "
$ws.Range("E49").Value = "  +1.73%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "81.06"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -3.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "73.02"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.32%  "
